# Applies the "Add files via upload" revision:
#   1. Refreshes the cached datetimeFigureOut field text (slide master,
#      all 5 slide layouts, notes master) from 31/08/2024 -> 09/09/2024.
#   2. Updates the student-details textbox on slide 1 (name, register
#      number, department, college).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Replace the first occurrence of $oldSub inside $tr's text with
# $newSub, touching only the characters that changed so existing run
# formatting (rPr) on the untouched parts is left alone.
function Replace-SubText($tr, $oldSub, $newSub) {
    $full = $tr.Text
    $idx = $full.IndexOf($oldSub)
    if ($idx -lt 0) {
        return $false
    }
    $tr.Characters($idx + 1, $oldSub.Length).Text = $newSub
    return $true
}

# Find a top-level shape in $shapes whose text exactly equals $targetText.
function Find-ShapeWithText($shapes, $targetText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $targetText) {
                    return $sh
                }
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. datetimeFigureOut fields: 31/08/2024 -> 09/09/2024
# ---------------------------------------------------------------------

$oldUS = "8/31/2024"
$newUS = "9/9/2024"
$oldDash = "31-08-2024"
$newDash = "09-09-2024"

# Slide master
$master = $p.SlideMaster
$dtShape = Find-ShapeWithText $master.Shapes $oldUS
if ($dtShape -ne $null) {
    $dtShape.TextFrame.TextRange.Text = $newUS
}

# Every slide layout under the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $dtShape = Find-ShapeWithText $layout.Shapes $oldUS
    if ($dtShape -ne $null) {
        $dtShape.TextFrame.TextRange.Text = $newUS
    }
}

# Notes master (different cached format: DD-MM-YYYY)
$notesMaster = $p.NotesMaster
$dtShape = Find-ShapeWithText $notesMaster.Shapes $oldDash
if ($dtShape -ne $null) {
    $dtShape.TextFrame.TextRange.Text = $newDash
}

# ---------------------------------------------------------------------
# 2. Slide 1 student-details textbox
# ---------------------------------------------------------------------

$slide1 = $p.Slides.Item(1)
$infoShape = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 13") {
        $infoShape = $sh
    }
}

if ($infoShape -ne $null) {
    $tr = $infoShape.TextFrame.TextRange

    # STUDENT NAME: SUBASHINI.J -> STUDENT NAME: HARISH K
    Replace-SubText $tr ": SUBASHINI.J" ": HARISH K"

    # REGISTER NO:  312210872 -> REGISTER NO:  312219553
    Replace-SubText $tr " 312210872" " 312219553"

    # DEPARTMENT: B.COM(General) III - B -> DEPARTMENT: B.COM(CA)III
    Replace-SubText $tr "(General) III - B" "(CA)III "

    # COLLEGE: BHAKTAVATSALAM MEMORIAL COLLEGE FOR WOMEN
    #   -> COLLEGE: SA COLLEGE OF ARTS AND SCIENCE
    Replace-SubText $tr "BHAKTAVATSALAM MEMORIAL COLLEGE FOR WOMEN" "SA COLLEGE OF ARTS AND SCIENCE "
}

"edit.ps1 completed"
